$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Range("I2:I209")

$col.Replace("POINT (650786.2464680928 6885812.586896558)", "POINT (5.846112318970617 52.65701914821788)") | Out-Null
$col.Replace("POINT (651894.4825726995 6887235.298555037)", "POINT (5.8560677732824304 52.66479009862622)") | Out-Null
$col.Replace("POINT (631822.1417841916 6878055.230200602)", "POINT (5.675754868098707 52.61462349689506)") | Out-Null
$col.Replace("POINT (631657.7821120179 6878055.884592394)", "POINT (5.674278400042641 52.61462707503078)") | Out-Null
$col.Replace("POINT (632503.7056526174 6888007.48048795)", "POINT (5.681877460499812 52.669007228992086)") | Out-Null
$col.Replace("POINT (632487.2493623364 6888007.547576849)", "POINT (5.681729631129018 52.66900759536806)") | Out-Null
$col.Replace("POINT (631717.705383289 6887750.637739368)", "POINT (5.674816699947214 52.667604574699354)") | Out-Null
$col.Replace("POINT (631716.0598059476 6887750.6442887895)", "POINT (5.674801917474444 52.667604610467244)") | Out-Null
$col.Replace("POINT (660779.3183178427 6885591.580651903)", "POINT (5.9358816107499655 52.65581187208736)") | Out-Null
$col.Replace("POINT (658899.1081426805 6885623.729513086)", "POINT (5.918991395372913 52.65598749161802)") | Out-Null
$col.Replace("POINT (662850.7493539307 6885819.082394737)", "POINT (5.95448959234714 52.65705463023741)") | Out-Null
$col.Replace("POINT (662867.2006905742 6885818.953272207)", "POINT (5.954637377218651 52.6570539248987)") | Out-Null
$col.Replace("POINT (648739.3188747525 6901165.176710194)", "POINT (5.827724455544781 52.74080273858885)") | Out-Null
$col.Replace("POINT (648739.3088234458 6901163.528410285)", "POINT (5.827724365252357 52.7407937520058)") | Out-Null
$col.Replace("POINT (651628.9345387137 6884589.821193228)", "POINT (5.8536823147064565 52.650339207109866)") | Out-Null
$col.Replace("POINT (651790.1350825747 6884588.779312476)", "POINT (5.855130403830043 52.65033351490178)") | Out-Null
$col.Replace("POINT (654088.4709180001 6894617.452217638)", "POINT (5.8757767059200665 52.70508964820048)") | Out-Null
$col.Replace("POINT (654088.4820535277 6894619.099167924)", "POINT (5.875776805952213 52.7050986348158)") | Out-Null
$col.Replace("POINT (630903.5278156842 6896765.268726382)", "POINT (5.667502818417548 52.71680766935835)") | Out-Null
$col.Replace("POINT (631262.4678542548 6895454.253963909)", "POINT (5.670727231644851 52.70965543315031)") | Out-Null
$col.Replace("POINT (627322.1392433867 6882912.21817465)", "POINT (5.635330657488889 52.64117283118716)") | Out-Null
$col.Replace("POINT (627125.3170264955 6882585.620353284)", "POINT (5.6335625734320125 52.63938808530209)") | Out-Null
$col.Replace("POINT (639798.0298312898 6902698.3963617105)", "POINT (5.747403489470051 52.74916109401226)") | Out-Null
$col.Replace("POINT (639814.5161055208 6902698.314131074)", "POINT (5.747551588191251 52.74916064577468)") | Out-Null
$col.Replace("POINT (623622.7808753886 6844055.484563485)", "POINT (5.602098755854809 52.428320073511365)") | Out-Null
$col.Replace("POINT (622533.6980640789 6845548.17729116)", "POINT (5.592315358504094 52.43651603250899)") | Out-Null
$col.Replace("POINT (639557.0575374895 6871422.563169483)", "POINT (5.745238798524349 52.57834177200539)") | Out-Null
$col.Replace("POINT (639422.5239744168 6871449.504573588)", "POINT (5.744030262964997 52.57848920677466)") | Out-Null
$col.Replace("POINT (629858.6037864573 6847117.042654197)", "POINT (5.658116106155565 52.445128583203314)") | Out-Null
$col.Replace("POINT (629887.9156513029 6847074.3626168175)", "POINT (5.658379419117534 52.44489430633676)") | Out-Null
$col.Replace("POINT (624811.4590703091 6843145.187948278)", "POINT (5.612776833758774 52.4233211370272)") | Out-Null
$col.Replace("POINT (624813.0955672998 6843145.182847288)", "POINT (5.6127915346613655 52.42332110901328)") | Out-Null
$col.Replace("POINT (619027.4140972253 6843571.916948595)", "POINT (5.560817873725216 52.42566461024887)") | Out-Null
$col.Replace("POINT (618998.8226262622 6841884.830314618)", "POINT (5.5605610321715995 52.41639888357573)") | Out-Null
$col.Replace("POINT (627750.3559306372 6848706.467726416)", "POINT (5.639177393439611 52.453852274568106)") | Out-Null
$col.Replace("POINT (627738.921064858 6848714.695758773)", "POINT (5.6390746722925975 52.45389743028845)") | Out-Null
$col.Replace("POINT (637329.1604652201 6855415.689658398)", "POINT (5.7252252586097025 52.49065726287612)") | Out-Null
$col.Replace("POINT (637329.16812207 6855417.328635846)", "POINT (5.725225327392356 52.49066625008681)") | Out-Null
$col.Replace("POINT (624443.40485739 6858660.661509362)", "POINT (5.609470546510277 52.50844720874757)") | Out-Null
$col.Replace("POINT (626182.8894762024 6858599.379804738)", "POINT (5.625096602705977 52.50811131038396)") | Out-Null
$col.Replace("POINT (632015.1153749432 6859817.7604676485)", "POINT (5.677488379358743 52.51478903579811)") | Out-Null
$col.Replace("POINT (631838.0280555513 6859823.388878135)", "POINT (5.6758975769024085 52.51481988173804)") | Out-Null
$col.Replace("POINT (612149.8555157806 6850606.275022324)", "POINT (5.499035713813825 52.4642772463945)") | Out-Null
$col.Replace("POINT (612149.8301782944 6850589.894837412)", "POINT (5.499035486203314 52.464187372634534)") | Out-Null
$col.Replace("POINT (615154.3498754905 6847379.718630455)", "POINT (5.526025545857608 52.44657042186665)") | Out-Null
$col.Replace("POINT (615046.5429271793 6847515.827626682)", "POINT (5.525057099563585 52.447317510792494)") | Out-Null
$col.Replace("POINT (608660.4861916584 6853207.604057717)", "POINT (5.467690175855856 52.47854771373215)") | Out-Null
$col.Replace("POINT (608660.467937963 6853191.218572761)", "POINT (5.467690011880121 52.478457840184284)") | Out-Null
$col.Replace("POINT (607868.7118403557 6852859.443142777)", "POINT (5.460577545842366 52.47663802921268)") | Out-Null
$col.Replace("POINT (607867.0733607284 6852859.444806904)", "POINT (5.460562827129448 52.47663803834072)") | Out-Null
$col.Replace("POINT (638862.7106414264 6858670.674769487)", "POINT (5.739001374232205 52.5085020933609)") | Out-Null
$col.Replace("POINT (638700.3865048984 6858671.462454493)", "POINT (5.73754319170396 52.508506410811684)") | Out-Null
$col.Replace("POINT (632198.5002912073 6842693.297169053)", "POINT (5.6791357540903125 52.42083934604499)") | Out-Null
$col.Replace("POINT (632196.8638929026 6842693.303768688)", "POINT (5.6791210540742325 52.42083938229131)") | Out-Null
$col.Replace("POINT (564591.5888895027 6918490.371498683)", "POINT (5.0718125358476565 52.835157406714806)") | Out-Null
$col.Replace("POINT (562939.7315095988 6918479.688106614)", "POINT (5.056973648532123 52.835099287289246)") | Out-Null
$col.Replace("POINT (564184.7322829454 6917104.407710172)", "POINT (5.068157680766503 52.82761688114269)") | Out-Null
$col.Replace("POINT (564052.1513856278 6917206.222421963)", "POINT (5.0669666863020755 52.82817086277342)") | Out-Null
$col.Replace("POINT (568027.1141880923 6908955.677622048)", "POINT (5.10267438469468 52.78325599507085)") | Out-Null
$col.Replace("POINT (568027.0492469635 6908972.1768434495)", "POINT (5.102673801318594 52.78334586116521)") | Out-Null
$col.Replace("POINT (558215.3713765166 6921344.754190781)", "POINT (5.014533999379798 52.850682922506486)") | Out-Null
$col.Replace("POINT (558215.3799004626 6921343.101738252)", "POINT (5.0145340759517065 52.85067393612504)") | Out-Null
$col.Replace("POINT (561800.3181996377 6915940.778253664)", "POINT (5.046738124619451 52.82128498012977)") | Out-Null
$col.Replace("POINT (561801.0962074423 6915775.645242449)", "POINT (5.046745113582472 52.8203863321109)") | Out-Null
$col.Replace("POINT (550894.2548993153 6907012.374623039)", "POINT (4.948767291096905 52.772670127312324)") | Out-Null
$col.Replace("POINT (552541.6501455882 6906537.24272937)", "POINT (4.963566094384033 52.770081519999074)") | Out-Null
$col.Replace("POINT (563630.5226015913 6908103.991852495)", "POINT (5.063179130492828 52.77861687599245)") | Out-Null
$col.Replace("POINT (563526.0252824926 6908228.9065566)", "POINT (5.06224041510387 52.77929731553227)") | Out-Null
$col.Replace("POINT (561658.5046929332 6916595.734906507)", "POINT (5.045464192213778 52.82484904720382)") | Out-Null
$col.Replace("POINT (561656.923506299 6916580.863777193)", "POINT (5.045449988172574 52.824768126478084)") | Out-Null
$col.Replace("POINT (565987.1137058334 6925667.158879593)", "POINT (5.084348748566437 52.87418266324956)") | Out-Null
$col.Replace("POINT (565988.7670570472 6925667.165810666)", "POINT (5.084363600873091 52.874182700921665)") | Out-Null
$col.Replace("POINT (567783.7039193974 6917913.961476444)", "POINT (5.100487793047877 52.832021528907575)") | Out-Null
$col.Replace("POINT (566171.4029992243 6917555.5803947)", "POINT (5.086004247455962 52.83007169224504)") | Out-Null
$col.Replace("POINT (565019.5471161997 6922295.814725481)", "POINT (5.075656950007723 52.85585467786636)") | Out-Null
$col.Replace("POINT (565016.1919129657 6922307.3690778185)", "POINT (5.075626809704259 52.85591750526866)") | Out-Null
$col.Replace("POINT (557976.9263298567 6907750.013986221)", "POINT (5.012392011081426 52.77668861789061)") | Out-Null
$col.Replace("POINT (557976.9177756517 6907751.663647342)", "POINT (5.0123919342376935 52.77669760444939)") | Out-Null
$col.Replace("POINT (608327.9250805442 6839529.476515237)", "POINT (5.46470272856568 52.403459662410285)") | Out-Null
$col.Replace("POINT (607092.4760638818 6840831.251975645)", "POINT (5.4536045012214975 52.410611484467815)") | Out-Null
$col.Replace("POINT (601854.9059225518 6826422.300072547)", "POINT (5.406554608125449 52.33138521651187)") | Out-Null
$col.Replace("POINT (601309.7514469102 6827943.004055552)", "POINT (5.401657402148699 52.33975340617269)") | Out-Null
$col.Replace("POINT (594010.0746288248 6825476.71179399)", "POINT (5.336083289600509 52.326180993749134)") | Out-Null
$col.Replace("POINT (594510.3798322884 6827041.543079308)", "POINT (5.340577607710467 52.334793003337516)") | Out-Null
$col.Replace("POINT (594030.0324088925 6824959.10479718)", "POINT (5.336262573389228 52.32333198578596)") | Out-Null
$col.Replace("POINT (595534.1803312879 6824337.937417271)", "POINT (5.349774564071872 52.319912718248375)") | Out-Null
$col.Replace("POINT (596825.2168961616 6820930.149456345)", "POINT (5.361372142857704 52.301149540964694)") | Out-Null
$col.Replace("POINT (596986.7843912605 6820933.469063375)", "POINT (5.362823528360347 52.30116782250991)") | Out-Null
$col.Replace("POINT (599588.2241068962 6819135.647732114)", "POINT (5.386192658933058 52.291265855867444)") | Out-Null
$col.Replace("POINT (599770.9693192092 6819243.336959521)", "POINT (5.387834287106262 52.2918590448272)") | Out-Null
$col.Replace("POINT (607978.5984596167 6835469.238615422)", "POINT (5.4615646741383905 52.38114564204505)") | Out-Null
$col.Replace("POINT (608074.9141276991 6835328.5344192125)", "POINT (5.462429892505776 52.38037216480185)") | Out-Null
$col.Replace("POINT (603116.5195819533 6839929.2472601915)", "POINT (5.417887876454393 52.405656086229605)") | Out-Null
$col.Replace("POINT (603208.0673500998 6839791.79712062)", "POINT (5.418710264047924 52.40490091890363)") | Out-Null
$col.Replace("POINT (604793.7235195715 6832903.7871513795)", "POINT (5.432954455771871 52.367040736821025)") | Out-Null
$col.Replace("POINT (604805.1717946685 6832915.220908703)", "POINT (5.433057297376837 52.36710360989914)") | Out-Null
$col.Replace("POINT (609936.203946678 6830774.4094128795)", "POINT (5.479150143431424 52.35532993541956)") | Out-Null
$col.Replace("POINT (609952.5439167006 6830774.388624065)", "POINT (5.479296927879559 52.355329821073425)") | Out-Null
$col.Replace("POINT (586532.1132754266 6825961.14678148)", "POINT (5.268907619822382 52.32884724959388)") | Out-Null
$col.Replace("POINT (586532.1398948386 6825944.81658633)", "POINT (5.268907858948628 52.328757373335584)") | Out-Null
$col.Replace("POINT (598055.834047391 6836712.729917548)", "POINT (5.372426964816193 52.38798075401587)") | Out-Null
$col.Replace("POINT (598055.8307490128 6836729.081956357)", "POINT (5.3724269351863585 52.388070629358296)") | Out-Null
$col.Replace("POINT (596726.3337596373 6836933.107475961)", "POINT (5.360483860528889 52.389191994761966)") | Out-Null
$col.Replace("POINT (596727.969006476 6836933.108075241)", "POINT (5.360498550201175 52.38919199805569)") | Out-Null
$col.Replace("POINT (605535.8116705157 6826227.763510216)", "POINT (5.439620747053444 52.33031459826143)") | Out-Null
$col.Replace("POINT (605535.8104855964 6826226.130432984)", "POINT (5.4396207364091325 52.33030561062583)") | Out-Null
$col.Replace("POINT (600115.7141596015 6827879.453226768)", "POINT (5.390931182698719 52.33940372803456)") | Out-Null
$col.Replace("POINT (600114.0807472114 6827879.453314338)", "POINT (5.3909165095055664 52.33940372851641)") | Out-Null
$col.Replace("POINT (593021.3229984895 6828461.497419653)", "POINT (5.327201182583226 52.342606228068114)") | Out-Null
$col.Replace("POINT (593022.9565280195 6828461.498768406)", "POINT (5.327215856828666 52.3426062354889)") | Out-Null

$ws.Range("J2").Value = 5.846112318970617
$ws.Range("K2").Value = 52.65701914821788
$ws.Range("J3").Value = 5.846112318970617
$ws.Range("K3").Value = 52.65701914821788
$ws.Range("J4").Value = 5.846112318970617
$ws.Range("K4").Value = 52.65701914821788
$ws.Range("J5").Value = 5.8560677732824304
$ws.Range("K5").Value = 52.66479009862622
$ws.Range("J6").Value = 5.8560677732824304
$ws.Range("K6").Value = 52.66479009862622
$ws.Range("J7").Value = 5.675754868098707
$ws.Range("K7").Value = 52.61462349689506
$ws.Range("J8").Value = 5.675754868098707
$ws.Range("K8").Value = 52.61462349689506
$ws.Range("J9").Value = 5.674278400042641
$ws.Range("K9").Value = 52.61462707503078
$ws.Range("J10").Value = 5.674278400042641
$ws.Range("K10").Value = 52.61462707503078
$ws.Range("J11").Value = 5.674278400042641
$ws.Range("K11").Value = 52.61462707503078
$ws.Range("J12").Value = 5.681877460499812
$ws.Range("K12").Value = 52.669007228992086
$ws.Range("J13").Value = 5.681877460499812
$ws.Range("K13").Value = 52.669007228992086
$ws.Range("J14").Value = 5.681877460499812
$ws.Range("K14").Value = 52.669007228992086
$ws.Range("J15").Value = 5.681729631129018
$ws.Range("K15").Value = 52.66900759536806
$ws.Range("J16").Value = 5.681729631129018
$ws.Range("K16").Value = 52.66900759536806
$ws.Range("J17").Value = 5.674816699947214
$ws.Range("K17").Value = 52.667604574699354
$ws.Range("J18").Value = 5.674816699947214
$ws.Range("K18").Value = 52.667604574699354
$ws.Range("J19").Value = 5.674801917474444
$ws.Range("K19").Value = 52.667604610467244
$ws.Range("J20").Value = 5.674801917474444
$ws.Range("K20").Value = 52.667604610467244
$ws.Range("J21").Value = 5.9358816107499655
$ws.Range("K21").Value = 52.65581187208736
$ws.Range("J22").Value = 5.9358816107499655
$ws.Range("K22").Value = 52.65581187208736
$ws.Range("J23").Value = 5.918991395372913
$ws.Range("K23").Value = 52.65598749161802
$ws.Range("J24").Value = 5.918991395372913
$ws.Range("K24").Value = 52.65598749161802
$ws.Range("J25").Value = 5.95448959234714
$ws.Range("K25").Value = 52.65705463023741
$ws.Range("J26").Value = 5.95448959234714
$ws.Range("K26").Value = 52.65705463023741
$ws.Range("J27").Value = 5.954637377218651
$ws.Range("K27").Value = 52.6570539248987
$ws.Range("J28").Value = 5.954637377218651
$ws.Range("K28").Value = 52.6570539248987
$ws.Range("J29").Value = 5.827724455544781
$ws.Range("K29").Value = 52.74080273858885
$ws.Range("J30").Value = 5.827724455544781
$ws.Range("K30").Value = 52.74080273858885
$ws.Range("J31").Value = 5.827724455544781
$ws.Range("K31").Value = 52.74080273858885
$ws.Range("J32").Value = 5.827724365252357
$ws.Range("K32").Value = 52.7407937520058
$ws.Range("J33").Value = 5.827724365252357
$ws.Range("K33").Value = 52.7407937520058
$ws.Range("J34").Value = 5.8536823147064565
$ws.Range("K34").Value = 52.650339207109866
$ws.Range("J35").Value = 5.8536823147064565
$ws.Range("K35").Value = 52.650339207109866
$ws.Range("J36").Value = 5.855130403830043
$ws.Range("K36").Value = 52.65033351490178
$ws.Range("J37").Value = 5.855130403830043
$ws.Range("K37").Value = 52.65033351490178
$ws.Range("J38").Value = 5.855130403830043
$ws.Range("K38").Value = 52.65033351490178
$ws.Range("J39").Value = 5.855130403830043
$ws.Range("K39").Value = 52.65033351490178
$ws.Range("J40").Value = 5.8757767059200665
$ws.Range("K40").Value = 52.70508964820048
$ws.Range("J41").Value = 5.8757767059200665
$ws.Range("K41").Value = 52.70508964820048
$ws.Range("J42").Value = 5.875776805952213
$ws.Range("K42").Value = 52.7050986348158
$ws.Range("J43").Value = 5.875776805952213
$ws.Range("K43").Value = 52.7050986348158
$ws.Range("J44").Value = 5.875776805952213
$ws.Range("K44").Value = 52.7050986348158
$ws.Range("J45").Value = 5.667502818417548
$ws.Range("K45").Value = 52.71680766935835
$ws.Range("J46").Value = 5.667502818417548
$ws.Range("K46").Value = 52.71680766935835
$ws.Range("J47").Value = 5.670727231644851
$ws.Range("K47").Value = 52.70965543315031
$ws.Range("J48").Value = 5.670727231644851
$ws.Range("K48").Value = 52.70965543315031
$ws.Range("J49").Value = 5.635330657488889
$ws.Range("K49").Value = 52.64117283118716
$ws.Range("J50").Value = 5.6335625734320125
$ws.Range("K50").Value = 52.63938808530209
$ws.Range("J51").Value = 5.6335625734320125
$ws.Range("K51").Value = 52.63938808530209
$ws.Range("J52").Value = 5.6335625734320125
$ws.Range("K52").Value = 52.63938808530209
$ws.Range("J53").Value = 5.747403489470051
$ws.Range("K53").Value = 52.74916109401226
$ws.Range("J54").Value = 5.747403489470051
$ws.Range("K54").Value = 52.74916109401226
$ws.Range("J55").Value = 5.747551588191251
$ws.Range("K55").Value = 52.74916064577468
$ws.Range("J56").Value = 5.747551588191251
$ws.Range("K56").Value = 52.74916064577468
$ws.Range("J57").Value = 5.602098755854809
$ws.Range("K57").Value = 52.428320073511365
$ws.Range("J58").Value = 5.602098755854809
$ws.Range("K58").Value = 52.428320073511365
$ws.Range("J59").Value = 5.592315358504094
$ws.Range("K59").Value = 52.43651603250899
$ws.Range("J60").Value = 5.745238798524349
$ws.Range("K60").Value = 52.57834177200539
$ws.Range("J61").Value = 5.745238798524349
$ws.Range("K61").Value = 52.57834177200539
$ws.Range("J62").Value = 5.744030262964997
$ws.Range("K62").Value = 52.57848920677466
$ws.Range("J63").Value = 5.744030262964997
$ws.Range("K63").Value = 52.57848920677466
$ws.Range("J64").Value = 5.744030262964997
$ws.Range("K64").Value = 52.57848920677466
$ws.Range("J65").Value = 5.658116106155565
$ws.Range("K65").Value = 52.445128583203314
$ws.Range("J66").Value = 5.658116106155565
$ws.Range("K66").Value = 52.445128583203314
$ws.Range("J67").Value = 5.658379419117534
$ws.Range("K67").Value = 52.44489430633676
$ws.Range("J68").Value = 5.658379419117534
$ws.Range("K68").Value = 52.44489430633676
$ws.Range("J69").Value = 5.612776833758774
$ws.Range("K69").Value = 52.4233211370272
$ws.Range("J70").Value = 5.612776833758774
$ws.Range("K70").Value = 52.4233211370272
$ws.Range("J71").Value = 5.6127915346613655
$ws.Range("K71").Value = 52.42332110901328
$ws.Range("J72").Value = 5.6127915346613655
$ws.Range("K72").Value = 52.42332110901328
$ws.Range("J73").Value = 5.6127915346613655
$ws.Range("K73").Value = 52.42332110901328
$ws.Range("J74").Value = 5.560817873725216
$ws.Range("K74").Value = 52.42566461024887
$ws.Range("J75").Value = 5.560817873725216
$ws.Range("K75").Value = 52.42566461024887
$ws.Range("J76").Value = 5.5605610321715995
$ws.Range("K76").Value = 52.41639888357573
$ws.Range("J77").Value = 5.5605610321715995
$ws.Range("K77").Value = 52.41639888357573
$ws.Range("J78").Value = 5.5605610321715995
$ws.Range("K78").Value = 52.41639888357573
$ws.Range("J79").Value = 5.639177393439611
$ws.Range("K79").Value = 52.453852274568106
$ws.Range("J80").Value = 5.639177393439611
$ws.Range("K80").Value = 52.453852274568106
$ws.Range("J81").Value = 5.6390746722925975
$ws.Range("K81").Value = 52.45389743028845
$ws.Range("J82").Value = 5.6390746722925975
$ws.Range("K82").Value = 52.45389743028845
$ws.Range("J83").Value = 5.7252252586097025
$ws.Range("K83").Value = 52.49065726287612
$ws.Range("J84").Value = 5.7252252586097025
$ws.Range("K84").Value = 52.49065726287612
$ws.Range("J85").Value = 5.7252252586097025
$ws.Range("K85").Value = 52.49065726287612
$ws.Range("J86").Value = 5.725225327392356
$ws.Range("K86").Value = 52.49066625008681
$ws.Range("J87").Value = 5.725225327392356
$ws.Range("K87").Value = 52.49066625008681
$ws.Range("J88").Value = 5.609470546510277
$ws.Range("K88").Value = 52.50844720874757
$ws.Range("J89").Value = 5.609470546510277
$ws.Range("K89").Value = 52.50844720874757
$ws.Range("J90").Value = 5.625096602705977
$ws.Range("K90").Value = 52.50811131038396
$ws.Range("J91").Value = 5.625096602705977
$ws.Range("K91").Value = 52.50811131038396
$ws.Range("J92").Value = 5.625096602705977
$ws.Range("K92").Value = 52.50811131038396
$ws.Range("J93").Value = 5.677488379358743
$ws.Range("K93").Value = 52.51478903579811
$ws.Range("J94").Value = 5.677488379358743
$ws.Range("K94").Value = 52.51478903579811
$ws.Range("J95").Value = 5.6758975769024085
$ws.Range("K95").Value = 52.51481988173804
$ws.Range("J96").Value = 5.6758975769024085
$ws.Range("K96").Value = 52.51481988173804
$ws.Range("J97").Value = 5.499035713813825
$ws.Range("K97").Value = 52.4642772463945
$ws.Range("J98").Value = 5.499035713813825
$ws.Range("K98").Value = 52.4642772463945
$ws.Range("J99").Value = 5.499035486203314
$ws.Range("K99").Value = 52.464187372634534
$ws.Range("J100").Value = 5.499035486203314
$ws.Range("K100").Value = 52.464187372634534
$ws.Range("J101").Value = 5.499035486203314
$ws.Range("K101").Value = 52.464187372634534
$ws.Range("J102").Value = 5.526025545857608
$ws.Range("K102").Value = 52.44657042186665
$ws.Range("J103").Value = 5.526025545857608
$ws.Range("K103").Value = 52.44657042186665
$ws.Range("J104").Value = 5.526025545857608
$ws.Range("K104").Value = 52.44657042186665
$ws.Range("J105").Value = 5.525057099563585
$ws.Range("K105").Value = 52.447317510792494
$ws.Range("J106").Value = 5.525057099563585
$ws.Range("K106").Value = 52.447317510792494
$ws.Range("J107").Value = 5.467690175855856
$ws.Range("K107").Value = 52.47854771373215
$ws.Range("J108").Value = 5.467690175855856
$ws.Range("K108").Value = 52.47854771373215
$ws.Range("J109").Value = 5.467690011880121
$ws.Range("K109").Value = 52.478457840184284
$ws.Range("J110").Value = 5.467690011880121
$ws.Range("K110").Value = 52.478457840184284
$ws.Range("J111").Value = 5.467690011880121
$ws.Range("K111").Value = 52.478457840184284
$ws.Range("J112").Value = 5.460577545842366
$ws.Range("K112").Value = 52.47663802921268
$ws.Range("J113").Value = 5.460577545842366
$ws.Range("K113").Value = 52.47663802921268
$ws.Range("J114").Value = 5.460562827129448
$ws.Range("K114").Value = 52.47663803834072
$ws.Range("J115").Value = 5.460562827129448
$ws.Range("K115").Value = 52.47663803834072
$ws.Range("J116").Value = 5.739001374232205
$ws.Range("K116").Value = 52.5085020933609
$ws.Range("J117").Value = 5.739001374232205
$ws.Range("K117").Value = 52.5085020933609
$ws.Range("J118").Value = 5.73754319170396
$ws.Range("K118").Value = 52.508506410811684
$ws.Range("J119").Value = 5.73754319170396
$ws.Range("K119").Value = 52.508506410811684
$ws.Range("J120").Value = 5.6791357540903125
$ws.Range("K120").Value = 52.42083934604499
$ws.Range("J121").Value = 5.6791357540903125
$ws.Range("K121").Value = 52.42083934604499
$ws.Range("J122").Value = 5.6791210540742325
$ws.Range("K122").Value = 52.42083938229131
$ws.Range("J123").Value = 5.6791210540742325
$ws.Range("K123").Value = 52.42083938229131
$ws.Range("J124").Value = 5.0718125358476565
$ws.Range("K124").Value = 52.835157406714806
$ws.Range("J125").Value = 5.0718125358476565
$ws.Range("K125").Value = 52.835157406714806
$ws.Range("J126").Value = 5.056973648532123
$ws.Range("K126").Value = 52.835099287289246
$ws.Range("J127").Value = 5.056973648532123
$ws.Range("K127").Value = 52.835099287289246
$ws.Range("J128").Value = 5.068157680766503
$ws.Range("K128").Value = 52.82761688114269
$ws.Range("J129").Value = 5.068157680766503
$ws.Range("K129").Value = 52.82761688114269
$ws.Range("J130").Value = 5.0669666863020755
$ws.Range("K130").Value = 52.82817086277342
$ws.Range("J131").Value = 5.0669666863020755
$ws.Range("K131").Value = 52.82817086277342
$ws.Range("J132").Value = 5.10267438469468
$ws.Range("K132").Value = 52.78325599507085
$ws.Range("J133").Value = 5.10267438469468
$ws.Range("K133").Value = 52.78325599507085
$ws.Range("J134").Value = 5.102673801318594
$ws.Range("K134").Value = 52.78334586116521
$ws.Range("J135").Value = 5.102673801318594
$ws.Range("K135").Value = 52.78334586116521
$ws.Range("J136").Value = 5.102673801318594
$ws.Range("K136").Value = 52.78334586116521
$ws.Range("J137").Value = 5.102673801318594
$ws.Range("K137").Value = 52.78334586116521
$ws.Range("J138").Value = 5.014533999379798
$ws.Range("K138").Value = 52.850682922506486
$ws.Range("J139").Value = 5.014533999379798
$ws.Range("K139").Value = 52.850682922506486
$ws.Range("J140").Value = 5.0145340759517065
$ws.Range("K140").Value = 52.85067393612504
$ws.Range("J141").Value = 5.0145340759517065
$ws.Range("K141").Value = 52.85067393612504
$ws.Range("J142").Value = 5.046738124619451
$ws.Range("K142").Value = 52.82128498012977
$ws.Range("J143").Value = 5.046738124619451
$ws.Range("K143").Value = 52.82128498012977
$ws.Range("J144").Value = 5.046738124619451
$ws.Range("K144").Value = 52.82128498012977
$ws.Range("J145").Value = 5.046745113582472
$ws.Range("K145").Value = 52.8203863321109
$ws.Range("J146").Value = 5.046745113582472
$ws.Range("K146").Value = 52.8203863321109
$ws.Range("J147").Value = 5.046745113582472
$ws.Range("K147").Value = 52.8203863321109
$ws.Range("J148").Value = 4.948767291096905
$ws.Range("K148").Value = 52.772670127312324
$ws.Range("J149").Value = 4.948767291096905
$ws.Range("K149").Value = 52.772670127312324
$ws.Range("J150").Value = 4.963566094384033
$ws.Range("K150").Value = 52.770081519999074
$ws.Range("J151").Value = 4.963566094384033
$ws.Range("K151").Value = 52.770081519999074
$ws.Range("J152").Value = 5.063179130492828
$ws.Range("K152").Value = 52.77861687599245
$ws.Range("J153").Value = 5.063179130492828
$ws.Range("K153").Value = 52.77861687599245
$ws.Range("J154").Value = 5.06224041510387
$ws.Range("K154").Value = 52.77929731553227
$ws.Range("J155").Value = 5.06224041510387
$ws.Range("K155").Value = 52.77929731553227
$ws.Range("J156").Value = 5.045464192213778
$ws.Range("K156").Value = 52.82484904720382
$ws.Range("J157").Value = 5.045464192213778
$ws.Range("K157").Value = 52.82484904720382
$ws.Range("J158").Value = 5.045464192213778
$ws.Range("K158").Value = 52.82484904720382
$ws.Range("J159").Value = 5.045449988172574
$ws.Range("K159").Value = 52.824768126478084
$ws.Range("J160").Value = 5.045449988172574
$ws.Range("K160").Value = 52.824768126478084
$ws.Range("J161").Value = 5.084348748566437
$ws.Range("K161").Value = 52.87418266324956
$ws.Range("J162").Value = 5.084348748566437
$ws.Range("K162").Value = 52.87418266324956
$ws.Range("J163").Value = 5.084348748566437
$ws.Range("K163").Value = 52.87418266324956
$ws.Range("J164").Value = 5.084363600873091
$ws.Range("K164").Value = 52.874182700921665
$ws.Range("J165").Value = 5.084363600873091
$ws.Range("K165").Value = 52.874182700921665
$ws.Range("J166").Value = 5.100487793047877
$ws.Range("K166").Value = 52.832021528907575
$ws.Range("J167").Value = 5.100487793047877
$ws.Range("K167").Value = 52.832021528907575
$ws.Range("J168").Value = 5.086004247455962
$ws.Range("K168").Value = 52.83007169224504
$ws.Range("J169").Value = 5.086004247455962
$ws.Range("K169").Value = 52.83007169224504
$ws.Range("J170").Value = 5.075656950007723
$ws.Range("K170").Value = 52.85585467786636
$ws.Range("J171").Value = 5.075656950007723
$ws.Range("K171").Value = 52.85585467786636
$ws.Range("J172").Value = 5.075626809704259
$ws.Range("K172").Value = 52.85591750526866
$ws.Range("J173").Value = 5.075626809704259
$ws.Range("K173").Value = 52.85591750526866
$ws.Range("J174").Value = 5.012392011081426
$ws.Range("K174").Value = 52.77668861789061
$ws.Range("J175").Value = 5.012392011081426
$ws.Range("K175").Value = 52.77668861789061
$ws.Range("J176").Value = 5.012392011081426
$ws.Range("K176").Value = 52.77668861789061
$ws.Range("J177").Value = 5.0123919342376935
$ws.Range("K177").Value = 52.77669760444939
$ws.Range("J178").Value = 5.46470272856568
$ws.Range("K178").Value = 52.403459662410285
$ws.Range("J179").Value = 5.4536045012214975
$ws.Range("K179").Value = 52.410611484467815
$ws.Range("J180").Value = 5.406554608125449
$ws.Range("K180").Value = 52.33138521651187
$ws.Range("J181").Value = 5.401657402148699
$ws.Range("K181").Value = 52.33975340617269
$ws.Range("J182").Value = 5.336083289600509
$ws.Range("K182").Value = 52.326180993749134
$ws.Range("J183").Value = 5.340577607710467
$ws.Range("K183").Value = 52.334793003337516
$ws.Range("J184").Value = 5.336262573389228
$ws.Range("K184").Value = 52.32333198578596
$ws.Range("J185").Value = 5.349774564071872
$ws.Range("K185").Value = 52.319912718248375
$ws.Range("J186").Value = 5.361372142857704
$ws.Range("K186").Value = 52.301149540964694
$ws.Range("J187").Value = 5.362823528360347
$ws.Range("K187").Value = 52.30116782250991
$ws.Range("J188").Value = 5.386192658933058
$ws.Range("K188").Value = 52.291265855867444
$ws.Range("J189").Value = 5.387834287106262
$ws.Range("K189").Value = 52.2918590448272
$ws.Range("J190").Value = 5.4615646741383905
$ws.Range("K190").Value = 52.38114564204505
$ws.Range("J191").Value = 5.462429892505776
$ws.Range("K191").Value = 52.38037216480185
$ws.Range("J192").Value = 5.417887876454393
$ws.Range("K192").Value = 52.405656086229605
$ws.Range("J193").Value = 5.418710264047924
$ws.Range("K193").Value = 52.40490091890363
$ws.Range("J194").Value = 5.432954455771871
$ws.Range("K194").Value = 52.367040736821025
$ws.Range("J195").Value = 5.433057297376837
$ws.Range("K195").Value = 52.36710360989914
$ws.Range("J196").Value = 5.479150143431424
$ws.Range("K196").Value = 52.35532993541956
$ws.Range("J197").Value = 5.479296927879559
$ws.Range("K197").Value = 52.355329821073425
$ws.Range("J198").Value = 5.268907619822382
$ws.Range("K198").Value = 52.32884724959388
$ws.Range("J199").Value = 5.268907858948628
$ws.Range("K199").Value = 52.328757373335584
$ws.Range("J200").Value = 5.372426964816193
$ws.Range("K200").Value = 52.38798075401587
$ws.Range("J201").Value = 5.3724269351863585
$ws.Range("K201").Value = 52.388070629358296
$ws.Range("J202").Value = 5.360483860528889
$ws.Range("K202").Value = 52.389191994761966
$ws.Range("J203").Value = 5.360498550201175
$ws.Range("K203").Value = 52.38919199805569
$ws.Range("J204").Value = 5.439620747053444
$ws.Range("K204").Value = 52.33031459826143
$ws.Range("J205").Value = 5.4396207364091325
$ws.Range("K205").Value = 52.33030561062583
$ws.Range("J206").Value = 5.390931182698719
$ws.Range("K206").Value = 52.33940372803456
$ws.Range("J207").Value = 5.3909165095055664
$ws.Range("K207").Value = 52.33940372851641
$ws.Range("J208").Value = 5.327201182583226
$ws.Range("K208").Value = 52.342606228068114
$ws.Range("J209").Value = 5.327215856828666
$ws.Range("K209").Value = 52.3426062354889

Write-Output "done"